$p = $ppt.ActivePresentation
Write-Host "Designs.Count: $($p.Designs.Count)"
for ($i=1; $i -le $p.Designs.Count; $i++) {
    $d = $p.Designs.Item($i)
    Write-Host "Design $i : $($d.Name)"
}
